$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "68.944.49"
$ws.Range("E2").Value = "  +2.36%  "
$ws.Range("D3").Value = "3.708.25"
$ws.Range("E3").Value = "  +0.94%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.06%  "
Set-TextValue "D5" "612.75"
$ws.Range("E5").Value = "  +8.10%  "
Set-TextValue "D6" "191.96"
$ws.Range("E6").Value = "  +10.09%  "
Set-TextValue "D7" "0.634"
$ws.Range("E8").Value = "  +0.11%  "
Set-TextValue "D9" "0.713"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("E10").Value = "  -2.20%  "
Set-TextValue "D11" "56.48"
$ws.Range("E11").Value = "  +9.30%  "
Set-TextValue "D12" "0.0000288"
$ws.Range("E12").Value = "  -2.83%  "
Set-TextValue "D13" "10.50"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "4.300.96"
$ws.Range("E14").Value = "  +1.14%  "
$ws.Range("D15").Value = "3.706.25"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("E16").Value = "  -0.14%  "
Set-TextValue "D17" "19.25"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "68.761.50"
$ws.Range("E20").Value = "  +2.23%  "
Set-TextValue "D21" "409.85"
$ws.Range("E21").Value = "  +1.15%  "
Set-TextValue "D22" "4.59"
$ws.Range("E22").Value = "  +2.07%  "
Set-TextValue "D23" "89.06"
$ws.Range("E23").Value = "  +1.76%  "
$ws.Range("E24").Value = "  -1.11%  "
Set-TextValue "D25" "12.75"
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("E27").Value = "  +1.51%  "
Set-TextValue "D28" "3.73"
$ws.Range("E28").Value = "  -1.25%  "
Set-TextValue "D29" "9.61"
$ws.Range("E29").Value = "  +1.68%  "
Set-TextValue "D30" "32.97"
$ws.Range("E30").Value = "  +0.84%  "
Set-TextValue "D31" "7.24"
$ws.Range("E31").Value = "  -7.95%  "
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("E33").Value = "  +3.79%  "
Set-TextValue "D34" "625.27"
$ws.Range("E34").Value = "  +5.18%  "
Set-TextValue "D35" "44.58"
$ws.Range("E35").Value = "  +1.52%  "
Set-TextValue "D36" "65.72"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("E37").Value = "  +3.38%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "0.0₃0809"
$ws.Range("E39").Value = "  -11.03%  "
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("E41").Value = "  +3.57%  "
$ws.Range("E42").Value = "  -0.95%  "
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("E44").Value = "  +1.51%  "
Set-TextValue "D45" "0.139"
$ws.Range("E45").Value = "  +4.02%  "
$ws.Range("D46").Value = "2.870.41"
$ws.Range("E46").Value = "  +5.71%  "
Set-TextValue "D47" "9.04"
$ws.Range("E47").Value = "  -3.32%  "
$ws.Range("E48").Value = "  +1.34%  "
Set-TextValue "D49" "3.13"
$ws.Range("E49").Value = "  +0.96%  "
Set-TextValue "D50" "141.24"
$ws.Range("E50").Value = "  -0.23%  "
